$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M5").Value = -58.432658
$ws.Range("N5").Value = -34.576446
$ws.Range("O5").Value = "Palermo"
$ws.Range("P5").Value = "Capital Sur"

$ws.Range("M72").Value = -58.452583
$ws.Range("N72").Value = -34.558668
$ws.Range("O72").Value = "Saavedra"
$ws.Range("P72").Value = "Capital Norte"

$ws.Range("M88").Value = -58.414185
$ws.Range("N88").Value = -34.64524
$ws.Range("O88").Value = "San Telmo"
$ws.Range("P88").Value = "Capital Sur"

$ws.Range("M89").Value = -58.404946
$ws.Range("N89").Value = -34.617251
$ws.Range("O89").Value = "Almagro"
$ws.Range("P89").Value = "Capital Sur"

$ws.Range("M90").Value = -58.460818
$ws.Range("N90").Value = -34.618934
$ws.Range("O90").Value = "Boedo"
$ws.Range("P90").Value = "Capital Sur"
